$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.318.70"
$ws.Range("E2").Value = "  -2.06%  "
$ws.Range("D3").Value = "2.636.24"
$ws.Range("E3").Value = "  -3.53%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.40"
$ws.Range("E5").Value = "  -0.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.01"
$ws.Range("E6").Value = "  -1.68%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.542"
$ws.Range("E8").Value = "  -0.67%  "
$ws.Range("D9").Value = "2.635.51"
$ws.Range("E9").Value = "  -3.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.145"
$ws.Range("E10").Value = "  +0.50%  "
$ws.Range("E11").Value = "  +1.29%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.363"
$ws.Range("E12").Value = "  -0.69%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.21"
$ws.Range("E13").Value = "  -2.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.87"
$ws.Range("E14").Value = "  -3.05%  "
$ws.Range("D15").Value = "3.120.40"
$ws.Range("E15").Value = "  -3.55%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000184"
$ws.Range("E16").Value = "  -2.89%  "
$ws.Range("D17").Value = "67.145.64"
$ws.Range("E17").Value = "  -2.31%  "
$ws.Range("D18").Value = "2.614.30"
$ws.Range("E18").Value = "  -5.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.77"
$ws.Range("E19").Value = "  -0.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.84"
$ws.Range("E20").Value = "  +1.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "363.08"
$ws.Range("E21").Value = "  -2.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.38"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.77"
$ws.Range("E23").Value = "  -3.94%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.86"
$ws.Range("E24").Value = "  +9.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.00"
$ws.Range("E25").Value = "  -5.84%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "70.54"
$ws.Range("E27").Value = "  -4.40%  "
$ws.Range("D28").Value = "2.773.06"
$ws.Range("E28").Value = "  -3.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000102"
$ws.Range("E29").Value = "  -3.86%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.993"
$ws.Range("E30").Value = "  -0.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "553.32"
$ws.Range("E31").Value = "  -6.98%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.00"
$ws.Range("E32").Value = "  -3.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.39"
$ws.Range("E33").Value = "  -4.28%  "
$ws.Range("E34").Value = "  -2.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.131"
$ws.Range("E35").Value = "  -1.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.54"
$ws.Range("E37").Value = "  -5.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "157.71"
$ws.Range("E38").Value = "  -2.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.35"
$ws.Range("E39").Value = "  -2.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.371"
$ws.Range("E40").Value = "  -2.86%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.24"
$ws.Range("E41").Value = "  -4.71%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.82"
$ws.Range("E42").Value = "  -5.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.91"
$ws.Range("E43").Value = "  -0.40%  "
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.51"
$ws.Range("E45").Value = "  -6.20%  "
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.08"
$ws.Range("E46").Value = "  -2.04%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₆0302"
$ws.Range("E47").Value = "  -3.73%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.592"
$ws.Range("E48").Value = "  -2.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "153.51"
$ws.Range("E49").Value = "  -2.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.87"
$ws.Range("E50").Value = "  -2.60%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.72"
$ws.Range("E51").Value = "  -4.40%  "
